$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 2-15 down to rows 3-16, preserving all
# their values (weekly roll: oldest displayed row drops off, a new row
# of data is prepended at row 2).
$src = $ws.Range("A2:T15")
$vals = $src.Value2()
$dst = $ws.Range("A3:T16")
$dst.Value2 = $vals

# Populate row 2 with the new weekly data point.
$ws.Range("D2").Value2 = 45092
$ws.Range("M2").Value2 = 220
$ws.Range("N2").Value2 = 16000
$ws.Range("O2").Value2 = 16000
$ws.Range("P2").Value2 = 16000
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("S2").Value2 = 889
$ws.Range("T2").Value2 = 18
